$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "265.90"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.50"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.118"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06103"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.568"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.480"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.358"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8205"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01336"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1591"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08040"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03462"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03213"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09223"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.728"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001628"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04640"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006280"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006139"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001068"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.725"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002715"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04596"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006993"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01058"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005971"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009905"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8029"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001125"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001901"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01241"

# --- Row 42/43 swap: CEJI and BKEXToken exchange positions with new data ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1117"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003472"
$ws.Range("E43").Value = "42CEJICEJI"
